$wb = $excel.ActiveWorkbook

$wsRun = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")

# Update execute column values on RUNMANAGER sheet: rows 2 and 3, column C -> "yes"
$wsRun.Range("C2").Value = "yes"
$wsRun.Range("C3").Value = "yes"

# Update selection on RUNMANAGER sheet
$wsRun.Range("A2").Select()

# Update selection on DATA sheet and make it the active sheet
$wsData.Activate()
$wsData.Range("A6:E6").Select()
